$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.38436040364177
$ws.Range("C2").Value = 0.128120134547257
$ws.Range("D2").Value = 8.47326156988115
$ws.Range("F2").Value = 0.0001

$ws.Range("B3").Value = 2.17735510967635
$ws.Range("C3").Value = 0.0151205215949747
